$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (B13/C13 = "849935 - Humberto Felipe da Silva", with no label
# in column A) is removed entirely; all rows below it shift up by one.
$ws.Rows.Item(13).Delete()

# After the shift, update the B/C text of the rows whose displayed value
# changed (row numbers below are the *new*, post-delete numbering).

# Row 10 "Objetivos:" -> now shows the docente text
$ws.Range("B10").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C10").Value = "849935 - Humberto Felipe da Silva"

# Row 13 "Programa resumido:" -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 "Programa:" -> "01/01/1996"
$ws.Range("B15").Value = "01/01/1996"
$ws.Range("C15").Value = "01/01/1996"

# Row 18 "Método:" -> docente text again
$ws.Range("B18").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C18").Value = "849935 - Humberto Felipe da Silva"

# Row 19 "Critério:" -> the evaluation-method paragraph
$ws.Range("B19").Value = "Avaliação da disciplina constará de uma avaliação escrita programa e de um seminário  a ser apresentado pelos alunos no final do semestre."
$ws.Range("C19").Value = "Avaliação da disciplina constará de uma avaliação escrita programa e de um seminário  a ser apresentado pelos alunos no final do semestre."

# Row 20 "Norma de recuperação:" -> the grading-criteria formula
$ws.Range("B20").Value = "A avaliação escrita programa = P1 Seminário = P2   MP =(P1+P2)/2."
$ws.Range("C20").Value = "A avaliação escrita programa = P1 Seminário = P2   MP =(P1+P2)/2."

# Row 21 "Bibliografia:" -> the recovery-norm paragraph
$ws.Range("B21").Value = "Na recuperação haverá uma aula de revisão e na semana seguinte uma avaliação escrita. A média final será a média simples entre MP e nota da recuperação."
$ws.Range("C21").Value = "Na recuperação haverá uma aula de revisão e na semana seguinte uma avaliação escrita. A média final será a média simples entre MP e nota da recuperação."
